$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark all touched ranges as Text format so numeric-looking strings
# (prices, "3"/"4" hour codes) are stored as text, matching the source data.
$ws.Range("B10:B19").NumberFormat = "@"
$ws.Range("B28:B34").NumberFormat = "@"
$ws.Range("C10:C19").NumberFormat = "@"
$ws.Range("C28:C34").NumberFormat = "@"
$ws.Range("D2:D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24:D28").NumberFormat = "@"
$ws.Range("D40:D42").NumberFormat = "@"
$ws.Range("D44:D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E10:E19").NumberFormat = "@"
$ws.Range("E28:E34").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Step 2: write the new values, row by row.
$ws.Range("D2").Value = '248.47'
$ws.Range("G2").Value = '4'
$ws.Range("D3").Value = '22.44'
$ws.Range("G3").Value = '4'
$ws.Range("D4").Value = '5.398'
$ws.Range("G4").Value = '4'
$ws.Range("D5").Value = '0.05706'
$ws.Range("G5").Value = '4'
$ws.Range("D6").Value = '3.413'
$ws.Range("G6").Value = '4'
$ws.Range("D7").Value = '6.323'
$ws.Range("G7").Value = '4'
$ws.Range("D8").Value = '0.8109'
$ws.Range("G8").Value = '4'
$ws.Range("D9").Value = '0.9223'
$ws.Range("G9").Value = '4'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '0.01133'
$ws.Range("E10").Value = '9OneONE'
$ws.Range("G10").Value = '4'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1423'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("G11").Value = '4'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.07456'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("G12").Value = '4'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '0.03118'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G13").Value = '4'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.03027'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("G14").Value = '4'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09352'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("G15").Value = '4'
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '3.722'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("G16").Value = '4'
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '0.001570'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("G17").Value = '4'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '0.04755'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("G18").Value = '4'
$ws.Range("B19").Value = 'UpBots'
$ws.Range("C19").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D19").Value = '0.01829'
$ws.Range("E19").Value = '18UpBotsUBXTBestin24h'
$ws.Range("G19").Value = '4'
$ws.Range("D20").Value = '0.006486'
$ws.Range("G20").Value = '4'
$ws.Range("G21").Value = '4'
$ws.Range("D22").Value = '0.001024'
$ws.Range("G22").Value = '4'
$ws.Range("G23").Value = '4'
$ws.Range("D24").Value = '3.701'
$ws.Range("G24").Value = '4'
$ws.Range("D25").Value = '2.164'
$ws.Range("G25").Value = '4'
$ws.Range("D26").Value = '0.3302'
$ws.Range("G26").Value = '4'
$ws.Range("D27").Value = '0.1307'
$ws.Range("G27").Value = '4'
$ws.Range("B28").Value = 'Spectre.aiUtilityToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("D28").Value = '--'
$ws.Range("E28").Value = '27Spectre.aiUtilityTokenSXUT'
$ws.Range("G28").Value = '4'
$ws.Range("B29").Value = 'LegolasExchange'
$ws.Range("C29").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("E29").Value = '28LegolasExchangeLGO'
$ws.Range("G29").Value = '4'
$ws.Range("B30").Value = 'BitZToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("E30").Value = '29BitZTokenBZ'
$ws.Range("G30").Value = '4'
$ws.Range("B31").Value = 'Birake'
$ws.Range("C31").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("E31").Value = '30BirakeBIR'
$ws.Range("G31").Value = '4'
$ws.Range("B32").Value = 'ZBToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("E32").Value = '31ZBTokenZB'
$ws.Range("G32").Value = '4'
$ws.Range("B33").Value = 'NashExchange'
$ws.Range("C33").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range("E33").Value = '32NashExchangeNEX'
$ws.Range("G33").Value = '4'
$ws.Range("B34").Value = 'AAXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("E34").Value = '33AAXTokenAAB'
$ws.Range("G34").Value = '4'
$ws.Range("G35").Value = '4'
$ws.Range("G36").Value = '4'
$ws.Range("G37").Value = '4'
$ws.Range("G38").Value = '4'
$ws.Range("G39").Value = '4'
$ws.Range("D40").Value = '0.03990'
$ws.Range("G40").Value = '4'
$ws.Range("D41").Value = '0.006884'
$ws.Range("G41").Value = '4'
$ws.Range("D42").Value = '0.1065'
$ws.Range("G42").Value = '4'
$ws.Range("G43").Value = '4'
$ws.Range("D44").Value = '0.007532'
$ws.Range("G44").Value = '4'
$ws.Range("D45").Value = '0.00005894'
$ws.Range("G45").Value = '4'
$ws.Range("G46").Value = '4'
$ws.Range("D47").Value = '0.5003'
$ws.Range("G47").Value = '4'
$ws.Range("G48").Value = '4'
$ws.Range("G49").Value = '4'
$ws.Range("D50").Value = '0.01011'
$ws.Range("G50").Value = '4'
$ws.Range("G51").Value = '4'

# Step 3: reset the style back to Normal (clears the temporary Text number
# format) while the underlying stored type remains Text, matching the original
# workbook which had no explicit style on these cells.
$ws.Range("B10:B19").Style = "Normal"
$ws.Range("B28:B34").Style = "Normal"
$ws.Range("C10:C19").Style = "Normal"
$ws.Range("C28:C34").Style = "Normal"
$ws.Range("D2:D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24:D28").Style = "Normal"
$ws.Range("D40:D42").Style = "Normal"
$ws.Range("D44:D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("E10:E19").Style = "Normal"
$ws.Range("E28:E34").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
